$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new project row (row 17) with Arun Kumar's details
$ws.Range("B17").Value = "Arun Kumar"
$ws.Range("C17").Value = "arunabi1819@gmail.com"
$ws.Range("D17").Value = 6374247005

# Turn the e-mail address in C17 into a mailto hyperlink (mirrors the other rows)
$ws.Hyperlinks.Add($ws.Range("C17"), "mailto:arunabi1819@gmail.com")

# Update the active selection shown when the sheet is opened
$ws.Range("G20").Select()

Write-Host "done"
